$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final "support Zone" (column C) ticker list for rows 2..23 (index 0..21 in column A)
$tickers = @(
    "NSE:ASMS",
    "NSE:BRNL",
    "NSE:CGCL",
    "NSE:FCSSOFT",
    "NSE:GANGESSECU",
    "NSE:GROBTEA",
    "NSE:GTLINFRA",
    "NSE:HMAAGRO",
    "NSE:JHS",
    "NSE:JMFINANCIL",
    "NSE:KHAITANLTD",
    "NSE:KOTAKBKETF",
    "NSE:KOTAKCONS",
    "NSE:LICNFNHGP",
    "NSE:LORDSCHLO",
    "NSE:MARINE",
    "NSE:MVGJL",
    "NSE:NBCC",
    "NSE:NDL",
    "NSE:NEXTMEDIA",
    "NSE:ORBTEXP",
    "NSE:SADBHAV"
)

# Copy the formatting of the existing "A2" index cell so new index cells (rows 13-23)
# pick up the same bordered/bold/centered style used by the rest of column A.
$ws.Range("A2").Copy()

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2

    # Column A: sequential index starting at 0
    $ws.Cells.Item($row, 1).Value = $i
    if ($row -gt 12) {
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    }

    # Column C: updated ticker symbol
    $ws.Cells.Item($row, 3).Value = $tickers[$i]

    # Ensure B/D/E/F exist for newly-added rows so the row shape matches the rest of the table
    if ($row -gt 12) {
        if ($ws.Cells.Item($row, 2).Value -eq $null) { $ws.Cells.Item($row, 2).Value = "" }
        $ws.Cells.Item($row, 4).Value = ""
        $ws.Cells.Item($row, 5).Value = ""
        $ws.Cells.Item($row, 6).Value = ""
    }
}

# Column E2 gains a new value (was empty before)
$ws.Range("E2").Value = "NSE:HINDUNILVR"
